# Update 想去人数 (F column) values across all sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 524
$ws.Range("F5").Value = 2322
$ws.Range("F7").Value = 8187
$ws.Range("F8").Value = 120
$ws.Range("F10").Value = 1625
$ws.Range("F11").Value = 1334
$ws.Range("F12").Value = 213
$ws.Range("F13").Value = 4539
$ws.Range("F14").Value = 6187
$ws.Range("F15").Value = 804
$ws.Range("F17").Value = 1268
$ws.Range("F19").Value = 490
$ws.Range("F20").Value = 6544
$ws.Range("F21").Value = 364
$ws.Range("F23").Value = 56
$ws.Range("F24").Value = 4414
$ws.Range("F25").Value = 325
$ws.Range("F26").Value = 725
$ws.Range("F27").Value = 2068
$ws.Range("F28").Value = 1202
$ws.Range("F29").Value = 361
$ws.Range("F31").Value = 81
$ws.Range("F32").Value = 57
$ws.Range("F34").Value = 93
$ws.Range("F40").Value = 172
$ws.Range("F41").Value = 1236
$ws.Range("F44").Value = 1204
$ws.Range("F47").Value = 203

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 699
$ws.Range("F7").Value = 407
$ws.Range("F8").Value = 420
$ws.Range("F10").Value = 223
$ws.Range("F15").Value = 203
$ws.Range("F17").Value = 115
$ws.Range("F22").Value = 130
$ws.Range("F26").Value = 186

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F4").Value = 471
$ws.Range("F6").Value = 1595
$ws.Range("F8").Value = 3164
$ws.Range("F9").Value = 1094
$ws.Range("F11").Value = 1541
$ws.Range("F12").Value = 1875
$ws.Range("F13").Value = 363
$ws.Range("F14").Value = 236

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 1595
$ws.Range("F5").Value = 524
$ws.Range("F7").Value = 3164
$ws.Range("F8").Value = 2322
$ws.Range("F9").Value = 120
$ws.Range("F10").Value = 1094
$ws.Range("F12").Value = 1625
$ws.Range("F13").Value = 1541
$ws.Range("F14").Value = 1334
$ws.Range("F15").Value = 699
$ws.Range("F16").Value = 213
$ws.Range("F17").Value = 1875
$ws.Range("F18").Value = 4539
$ws.Range("F19").Value = 407
$ws.Range("F20").Value = 420
$ws.Range("F21").Value = 805
$ws.Range("F23").Value = 1268
$ws.Range("F25").Value = 490
$ws.Range("F26").Value = 6544
$ws.Range("F27").Value = 364
$ws.Range("F28").Value = 236
$ws.Range("F30").Value = 325
$ws.Range("F31").Value = 2068
$ws.Range("F32").Value = 1202
$ws.Range("F33").Value = 361
$ws.Range("F34").Value = 81
$ws.Range("F36").Value = 203
$ws.Range("F38").Value = 93
$ws.Range("F43").Value = 1236
$ws.Range("F44").Value = 130
$ws.Range("F47").Value = 1204
$ws.Range("F49").Value = 203
